$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano statistic (column C) and p-value (column D) for rows 2-11
$ws.Range("C2").Value = 0.4108978844957375
$ws.Range("D2").Value = 0.6851222294128634

$ws.Range("C3").Value = -0.04902813066291979
$ws.Range("D3").Value = 0.9613392686794529

$ws.Range("C4").Value = 0.8171163532916009
$ws.Range("D4").Value = 0.4226229798549848

$ws.Range("C5").Value = 0.4785386810923075
$ws.Range("D5").Value = 0.6369883361377542

$ws.Range("C6").Value = -0.5994576667036519
$ws.Range("D6").Value = 0.5549917442185048

$ws.Range("C7").Value = 0.2781818926557545
$ws.Range("D7").Value = 0.7834704883977843

$ws.Range("C8").Value = 0.1350207046742724
$ws.Range("D8").Value = 0.8938232320861812

$ws.Range("C9").Value = 1.170866014588386
$ws.Range("D9").Value = 0.2541792372173757

$ws.Range("C10").Value = 0.7460440740220322
$ws.Range("D10").Value = 0.4635386211443322

$ws.Range("C11").Value = -0.2402110480742554
$ws.Range("D11").Value = 0.812390821364205
